# events_criteria_tables.xlsx
#
# Adds "short label" lookup columns (column C) to the J. Cause,
# K. Impact and L. Action Taken tables -- these were pasted in from the
# Power BI (.pbix) data-model snippets referenced in the commit message,
# which is why the shared-string table below is populated in the exact
# "paste" order (J. Cause, then L. Action Taken, then K. Impact; and with
# a couple of rows out of natural row-order inside each) rather than a
# simple top-to-bottom fill. It also fills in the one placeholder " "
# criteria cell on "Alert Criteria" and restores the sheet selections /
# active tab to match the author's last save.

$wb = $excel.ActiveWorkbook

$wsAlert  = $wb.Worksheets.Item("Alert Criteria")
$wsCause  = $wb.Worksheets.Item("J. Cause")
$wsImpact = $wb.Worksheets.Item("K. Impact")
$wsAction = $wb.Worksheets.Item("L. Action Taken")

# ---------------------------------------------------------------------
# Alert Criteria: the placeholder " " text for row 14 (System Report /
# Physical threat to Facility) gets filled in with the real criteria text.
# ---------------------------------------------------------------------
$wsAlert.Range("B15").Value = " Damage or destruction of its Facility that results from actual or suspected intentional human action."

# ---------------------------------------------------------------------
# J. Cause -- column C ("Cause Short"). Written in the same order the
# values were originally pasted in, so newly-introduced shared strings
# land at the same table offsets as the source workbook.
# ---------------------------------------------------------------------
$wsCause.Range("C1").Value = "Cause Short"
$wsCause.Range("C8").Value = "Cyber event (information)"
$wsCause.Range("C9").Value = "Cyber event (operational)"
$wsCause.Range("C10").Value = "Fuel supply"
$wsCause.Range("C12").Value = "Transmission equipment failure "
$wsCause.Range("C13").Value = "Failure at high voltage"
$wsCause.Range("C11").Value = "Generator loss or failure "

# L. Action Taken -- column C ("Action Taken Short").
$wsAction.Range("C1").Value = "Action Taken Short"
$wsAction.Range("C3").Value = "Shed Firm Load "
$wsAction.Range("C4").Value = "Public Appeal "
$wsAction.Range("C5").Value = "Implemented warning etc."

# K. Impact -- column C ("Impact Short").
$wsImpact.Range("C1").Value = "Impact Short"
$wsImpact.Range("C3").Value = "Control center loss"
$wsImpact.Range("C5").Value = "Facility Damage"
$wsImpact.Range("C4").Value = "Control Center Communication"
$wsImpact.Range("C6").Value = "Electrical Separation"
$wsImpact.Range("C7").Value = "Complete Shutdown"
$wsImpact.Range("C9").Value = "Major Interruption"
$wsImpact.Range("C8").Value = "Major Interruption of +3 BES"
$wsImpact.Range("C10").Value = "Uncontrolled 200MW loss"
$wsImpact.Range("C11").Value = "Service Loss +50K customers"
$wsImpact.Range("C12").Value = "Voltage Reduction"
$wsImpact.Range("C13").Value = "Voltage Deviation"
$wsImpact.Range("C14").Value = "Inadequate resources "
$wsImpact.Range("C15").Value = "Cap. Loss + 1,400MW"
$wsImpact.Range("C16").Value = "Cap. Loss + 2,000MW"
$wsImpact.Range("C17").Value = "Loss Nuclear Generation"

# ---------------------------------------------------------------------
# Remaining column-C cells: these all re-use strings that already exist
# in the shared-string table (ID labels shared with column B, or with
# other sheets), so write order has no effect on the table layout.
# ---------------------------------------------------------------------
$wsCause.Range("C2").Value = "Unknown"
$wsCause.Range("C3").Value = "Physical attack"
$wsCause.Range("C4").Value = "Threat of physical attack"
$wsCause.Range("C5").Value = "Vandalism"
$wsCause.Range("C6").Value = "Theft"
$wsCause.Range("C7").Value = "Suspicious activity"
$wsCause.Range("C14").Value = "Weather or natural disaster"
$wsCause.Range("C15").Value = "Operator action(s)"
$wsCause.Range("C16").Value = "Other"

$wsAction.Range("C2").Value = "None"
$wsAction.Range("C6").Value = "Voltage reduction"
$wsAction.Range("C7").Value = "Shed Interruptible Load"
$wsAction.Range("C8").Value = "Repaired or restored"
$wsAction.Range("C9").Value = "Mitigation implemented"
$wsAction.Range("C10").Value = "Other"

$wsImpact.Range("C2").Value = "None"
$wsImpact.Range("C18").Value = "Other"

# ---------------------------------------------------------------------
# New column widths (manually sized, not auto "best fit").
# ---------------------------------------------------------------------
$wsCause.Columns.Item(3).ColumnWidth = 20.1667
$wsImpact.Columns.Item(3).ColumnWidth = 12.42
$wsAction.Columns.Item(3).ColumnWidth = 16.92

# ---------------------------------------------------------------------
# View state: restore each sheet's selection, and land on "K. Impact"
# as the active tab/sheet, matching the author's last save.
# ---------------------------------------------------------------------
$wsAlert.Range("C1:C1048576").Select()
$wsCause.Range("F21").Select()
$wsAction.Range("F40").Select()

$wsImpact.Activate()
$wsImpact.Range("D18").Select()
